# Generate Report for Handback
#
# - Overview sheet + zh-cn/de-de per-file sheets move from "Ready for
#   handoff" to "Handed back: in sync with en-US" for every cell that
#   carried the old status text.
# - The zh-cn and de-de sheets each gain a "Latest Target File" (F) and
#   "Latest Handback File" (G) hyperlink pair per row: F mirrors the
#   existing "Source File Name" (A) hyperlink/text, G mirrors the
#   existing "Latest Handoff File" (D) hyperlink/text -- the handback
#   confirms the target + handback files are the ones already handed off.
# - "Latest Handback DateTime" (H) moves off the zero-date placeholder:
#   zh-cn gets 2016-03-20 04:44:33, de-de gets 2016-03-20 04:44:47.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: refresh the status text wherever it appears ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in 2, 3) {
    if ($wsOverview.Cells.Item($r, 2).Value2 -eq $oldStatus) {
        $wsOverview.Cells.Item($r, 2).Value = $newStatus
    }
    if ($wsOverview.Cells.Item($r, 3).Value2 -eq $oldStatus) {
        $wsOverview.Cells.Item($r, 3).Value = $newStatus
    }
}

function Get-HyperlinkAddress($ws, $addr) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            return $hl.Address
        }
    }
    return $null
}

function Update-LangSheet($SheetName, $HandbackDateTime) {
    $ws = $wb.Worksheets.Item($SheetName)

    # Status column (C) -> new handed-back text
    if ($ws.Range("C2").Value2 -eq $oldStatus) {
        $ws.Range("C2").Value = $newStatus
    }
    if ($ws.Range("C3").Value2 -eq $oldStatus) {
        $ws.Range("C3").Value = $newStatus
    }

    foreach ($row in 2, 3) {
        $aAddr = "`$A`$$row"
        $dAddr = "`$D`$$row"

        $sourceUrl = Get-HyperlinkAddress $ws $aAddr
        $handoffUrl = Get-HyperlinkAddress $ws $dAddr

        $sourceText = $ws.Range("A$row").Value2
        $handoffText = $ws.Range("D$row").Value2

        # F = Latest Target File (mirrors Source File Name / column A)
        $ws.Hyperlinks.Add($ws.Range("F$row"), $sourceUrl, $null, $null, $sourceText) | Out-Null
        # G = Latest Handback File (mirrors Latest Handoff File / column D)
        $ws.Hyperlinks.Add($ws.Range("G$row"), $handoffUrl, $null, $null, $handoffText) | Out-Null
    }

    # Latest Handback DateTime (H)
    $ws.Range("H2").Value = $HandbackDateTime
    $ws.Range("H3").Value = $HandbackDateTime
}

Update-LangSheet "zh-cn" "2016-03-20 04:44:33"
Update-LangSheet "de-de" "2016-03-20 04:44:47"
